$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.938.19'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.363.47'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.83%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -3.10%  '
$ws.Range('D9').Value = '2.361.08'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.345'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.89'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '2.773.20'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '60.886.54'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '2.373.81'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.70'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.22'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.61'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.91'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -15.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.23'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '2.472.28'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.06'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('D31').Value = '0.0₃0873'
$ws.Range('E31').Value = '  -7.73%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '494.42'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -8.81%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.37'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.04%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.65'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.87'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.376'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.49'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.29'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '145.02'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.09%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.93'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '146.22'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.97%  '
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.01'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -9.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0517'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.571'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.10'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0904'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.75%  '
